# ----------------------------------------------------------------------
# "added TFIDF 1,2 KT"
#
# 1. TFIDF12 worksheet: the B1/J1 "one seed only" header labels are
#    cleared, and the raw numeric results in B3:G14 / I3:N14 are
#    replaced with "mean±std" text summaries (the underlying figures
#    for the TFIDF 1+2 KT experiment). The sheet's selection moves
#    to G26.
# 2. BOW2 becomes the active sheet/tab (previously TFIDF1 was active).
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsTFIDF12 = $wb.Worksheets.Item("TFIDF12")

# Clear the "one seed only" note cells (keep their existing style).
$wsTFIDF12.Cells.Item(1, 2).Value = ""
$wsTFIDF12.Cells.Item(1, 10).Value = ""

# Replace the raw numeric scores with their "mean±std" text summaries.
$wsTFIDF12.Cells.Item(3, 2).Value = "0.7416±0.0022"
$wsTFIDF12.Cells.Item(3, 3).Value = "0.6984±0.0059"
$wsTFIDF12.Cells.Item(3, 4).Value = "0.5981±0.0056"
$wsTFIDF12.Cells.Item(3, 5).Value = "0.4725±0.0048"
$wsTFIDF12.Cells.Item(3, 6).Value = "0.8207±0.002"
$wsTFIDF12.Cells.Item(3, 7).Value = "0.6444±0.0037"
$wsTFIDF12.Cells.Item(3, 9).Value = "0.7474±0.0023"
$wsTFIDF12.Cells.Item(3, 10).Value = "0.7087±0.0061"
$wsTFIDF12.Cells.Item(3, 11).Value = "0.6045±0.0039"
$wsTFIDF12.Cells.Item(3, 12).Value = "0.4849±0.0047"
$wsTFIDF12.Cells.Item(3, 13).Value = "0.8276±0.002"
$wsTFIDF12.Cells.Item(3, 14).Value = "0.6524±0.0035"
$wsTFIDF12.Cells.Item(4, 2).Value = "0.7387±0.002"
$wsTFIDF12.Cells.Item(4, 3).Value = "0.6871±0.0049"
$wsTFIDF12.Cells.Item(4, 4).Value = "0.6018±0.0038"
$wsTFIDF12.Cells.Item(4, 5).Value = "0.4695±0.0043"
$wsTFIDF12.Cells.Item(4, 6).Value = "0.8147±0.0015"
$wsTFIDF12.Cells.Item(4, 7).Value = "0.6417±0.0038"
$wsTFIDF12.Cells.Item(4, 9).Value = "0.7299±0.0026"
$wsTFIDF12.Cells.Item(4, 10).Value = "0.6621±0.0044"
$wsTFIDF12.Cells.Item(4, 11).Value = "0.6157±0.0024"
$wsTFIDF12.Cells.Item(4, 12).Value = "0.4647±0.0045"
$wsTFIDF12.Cells.Item(4, 13).Value = "0.7983±0.0023"
$wsTFIDF12.Cells.Item(4, 14).Value = "0.6381±0.003"
$wsTFIDF12.Cells.Item(5, 2).Value = "0.7165±0.0022"
$wsTFIDF12.Cells.Item(5, 3).Value = "0.6776±0.0054"
$wsTFIDF12.Cells.Item(5, 4).Value = "0.5327±0.0041"
$wsTFIDF12.Cells.Item(5, 5).Value = "0.4033±0.0057"
$wsTFIDF12.Cells.Item(5, 6).Value = "0.8087±0.0018"
$wsTFIDF12.Cells.Item(5, 7).Value = "0.5965±0.0041"
$wsTFIDF12.Cells.Item(5, 9).Value = "0.7211±0.0018"
$wsTFIDF12.Cells.Item(5, 10).Value = "0.6833±0.0056"
$wsTFIDF12.Cells.Item(5, 11).Value = "0.5441±0.0018"
$wsTFIDF12.Cells.Item(5, 12).Value = "0.416±0.0037"
$wsTFIDF12.Cells.Item(5, 13).Value = "0.8132±0.0034"
$wsTFIDF12.Cells.Item(5, 14).Value = "0.6058±0.0024"
$wsTFIDF12.Cells.Item(6, 2).Value = "0.7237±0.0022"
$wsTFIDF12.Cells.Item(6, 3).Value = "0.6861±0.0048"
$wsTFIDF12.Cells.Item(6, 4).Value = "0.5505±0.0039"
$wsTFIDF12.Cells.Item(6, 5).Value = "0.4231±0.0054"
$wsTFIDF12.Cells.Item(6, 6).Value = "0.8143±0.0029"
$wsTFIDF12.Cells.Item(6, 7).Value = "0.6109±0.0038"
$wsTFIDF12.Cells.Item(6, 9).Value = "0.7264±0.0026"
$wsTFIDF12.Cells.Item(6, 10).Value = "0.691±0.0062"
$wsTFIDF12.Cells.Item(6, 11).Value = "0.5543±0.0036"
$wsTFIDF12.Cells.Item(6, 12).Value = "0.4295±0.0062"
$wsTFIDF12.Cells.Item(6, 13).Value = "0.8214±0.0036"
$wsTFIDF12.Cells.Item(6, 14).Value = "0.6151±0.0038"
$wsTFIDF12.Cells.Item(7, 2).Value = "0.7107±0.003"
$wsTFIDF12.Cells.Item(7, 3).Value = "0.706±0.0068"
$wsTFIDF12.Cells.Item(7, 4).Value = "0.5029±0.0048"
$wsTFIDF12.Cells.Item(7, 5).Value = "0.383±0.0083"
$wsTFIDF12.Cells.Item(7, 6).Value = "0.8053±0.0024"
$wsTFIDF12.Cells.Item(7, 7).Value = "0.5874±0.0052"
$wsTFIDF12.Cells.Item(7, 9).Value = "0.7135±0.0015"
$wsTFIDF12.Cells.Item(7, 10).Value = "0.7179±0.0057"
$wsTFIDF12.Cells.Item(7, 11).Value = "0.5048±0.0031"
$wsTFIDF12.Cells.Item(7, 12).Value = "0.3908±0.0042"
$wsTFIDF12.Cells.Item(7, 13).Value = "0.8091±0.0031"
$wsTFIDF12.Cells.Item(7, 14).Value = "0.5928±0.0034"
$wsTFIDF12.Cells.Item(8, 2).Value = "0.7408±0.0019"
$wsTFIDF12.Cells.Item(8, 3).Value = "0.7037±0.0041"
$wsTFIDF12.Cells.Item(8, 4).Value = "0.5898±0.0028"
$wsTFIDF12.Cells.Item(8, 5).Value = "0.468±0.0042"
$wsTFIDF12.Cells.Item(8, 6).Value = "0.838±0.002"
$wsTFIDF12.Cells.Item(8, 7).Value = "0.6417±0.0027"
$wsTFIDF12.Cells.Item(8, 9).Value = "0.7426±0.0021"
$wsTFIDF12.Cells.Item(8, 10).Value = "0.7095±0.0042"
$wsTFIDF12.Cells.Item(8, 11).Value = "0.5883±0.004"
$wsTFIDF12.Cells.Item(8, 12).Value = "0.4708±0.0051"
$wsTFIDF12.Cells.Item(8, 13).Value = "0.8433±0.0026"
$wsTFIDF12.Cells.Item(8, 14).Value = "0.6433±0.0034"
$wsTFIDF12.Cells.Item(9, 2).Value = "0.7199±0.0037"
$wsTFIDF12.Cells.Item(9, 3).Value = "0.6478±0.006"
$wsTFIDF12.Cells.Item(9, 4).Value = "0.619±0.0059"
$wsTFIDF12.Cells.Item(9, 5).Value = "0.4533±0.0075"
$wsTFIDF12.Cells.Item(9, 6).Value = "0.82±0.004"
$wsTFIDF12.Cells.Item(9, 7).Value = "0.633±0.0055"
$wsTFIDF12.Cells.Item(9, 9).Value = "0.7273±0.004"
$wsTFIDF12.Cells.Item(9, 10).Value = "0.6596±0.0067"
$wsTFIDF12.Cells.Item(9, 11).Value = "0.6242±0.0043"
$wsTFIDF12.Cells.Item(9, 12).Value = "0.465±0.006"
$wsTFIDF12.Cells.Item(9, 13).Value = "0.8286±0.0026"
$wsTFIDF12.Cells.Item(9, 14).Value = "0.6414±0.0039"
$wsTFIDF12.Cells.Item(10, 2).Value = "0.725±0.0025"
$wsTFIDF12.Cells.Item(10, 3).Value = "0.6687±0.0043"
$wsTFIDF12.Cells.Item(10, 4).Value = "0.5822±0.0047"
$wsTFIDF12.Cells.Item(10, 5).Value = "0.4383±0.006"
$wsTFIDF12.Cells.Item(10, 6).Value = "0.8212±0.0016"
$wsTFIDF12.Cells.Item(10, 7).Value = "0.6225±0.0042"
$wsTFIDF12.Cells.Item(10, 9).Value = "0.7285±0.0025"
$wsTFIDF12.Cells.Item(10, 10).Value = "0.6744±0.0042"
$wsTFIDF12.Cells.Item(10, 11).Value = "0.5879±0.0042"
$wsTFIDF12.Cells.Item(10, 12).Value = "0.4465±0.0058"
$wsTFIDF12.Cells.Item(10, 13).Value = "0.8247±0.0027"
$wsTFIDF12.Cells.Item(10, 14).Value = "0.6282±0.0039"
$wsTFIDF12.Cells.Item(11, 2).Value = "0.6305±0.0046"
$wsTFIDF12.Cells.Item(11, 3).Value = "0.5103±0.0064"
$wsTFIDF12.Cells.Item(11, 4).Value = "0.4672±0.0045"
$wsTFIDF12.Cells.Item(11, 5).Value = "0.237±0.009"
$wsTFIDF12.Cells.Item(11, 6).Value = "0.6004±0.0034"
$wsTFIDF12.Cells.Item(11, 7).Value = "0.4878±0.0051"
$wsTFIDF12.Cells.Item(11, 9).Value = "0.6341±0.0051"
$wsTFIDF12.Cells.Item(11, 10).Value = "0.5178±0.0094"
$wsTFIDF12.Cells.Item(11, 11).Value = "0.4749±0.006"
$wsTFIDF12.Cells.Item(11, 12).Value = "0.2477±0.0098"
$wsTFIDF12.Cells.Item(11, 13).Value = "0.6062±0.0045"
$wsTFIDF12.Cells.Item(11, 14).Value = "0.4954±0.0075"
$wsTFIDF12.Cells.Item(12, 2).Value = "0.6346±0.0028"
$wsTFIDF12.Cells.Item(12, 3).Value = "0.529±0.0037"
$wsTFIDF12.Cells.Item(12, 4).Value = "0.5177±0.0053"
$wsTFIDF12.Cells.Item(12, 5).Value = "0.2941±0.0072"
$wsTFIDF12.Cells.Item(12, 6).Value = "0.6416±0.0043"
$wsTFIDF12.Cells.Item(12, 7).Value = "0.5233±0.0044"
$wsTFIDF12.Cells.Item(12, 9).Value = "0.6402±0.0037"
$wsTFIDF12.Cells.Item(12, 10).Value = "0.5357±0.0047"
$wsTFIDF12.Cells.Item(12, 11).Value = "0.5228±0.0038"
$wsTFIDF12.Cells.Item(12, 12).Value = "0.3034±0.0066"
$wsTFIDF12.Cells.Item(12, 13).Value = "0.6448±0.0032"
$wsTFIDF12.Cells.Item(12, 14).Value = "0.5292±0.0041"
$wsTFIDF12.Cells.Item(13, 2).Value = "0.7395±0.0025"
$wsTFIDF12.Cells.Item(13, 3).Value = "0.6871±0.0035"
$wsTFIDF12.Cells.Item(13, 4).Value = "0.6048±0.0038"
$wsTFIDF12.Cells.Item(13, 5).Value = "0.4721±0.0056"
$wsTFIDF12.Cells.Item(13, 6).Value = "0.8377±0.0015"
$wsTFIDF12.Cells.Item(13, 7).Value = "0.6433±0.0033"
$wsTFIDF12.Cells.Item(13, 9).Value = "0.7439±0.002"
$wsTFIDF12.Cells.Item(13, 10).Value = "0.6944±0.0043"
$wsTFIDF12.Cells.Item(13, 11).Value = "0.6093±0.0038"
$wsTFIDF12.Cells.Item(13, 12).Value = "0.4811±0.0046"
$wsTFIDF12.Cells.Item(13, 13).Value = "0.8444±0.0026"
$wsTFIDF12.Cells.Item(13, 14).Value = "0.6491±0.0034"
$wsTFIDF12.Cells.Item(14, 2).Value = "0.7237±0.0035"
$wsTFIDF12.Cells.Item(14, 3).Value = "0.6925±0.0055"
$wsTFIDF12.Cells.Item(14, 4).Value = "0.5452±0.0074"
$wsTFIDF12.Cells.Item(14, 5).Value = "0.4213±0.0093"
$wsTFIDF12.Cells.Item(14, 6).Value = "0.8244±0.0017"
$wsTFIDF12.Cells.Item(14, 7).Value = "0.6101±0.0056"
$wsTFIDF12.Cells.Item(14, 9).Value = "0.7254±0.002"
$wsTFIDF12.Cells.Item(14, 10).Value = "0.7012±0.0042"
$wsTFIDF12.Cells.Item(14, 11).Value = "0.5441±0.0032"
$wsTFIDF12.Cells.Item(14, 12).Value = "0.4247±0.0049"
$wsTFIDF12.Cells.Item(14, 13).Value = "0.827±0.0028"
$wsTFIDF12.Cells.Item(14, 14).Value = "0.6128±0.0029"

# Update the sheet's stored selection.
$wsTFIDF12.Activate()
$wsTFIDF12.Range("G26").Select()

# BOW2 is now the active/selected tab in the workbook.
$wsBOW2 = $wb.Worksheets.Item("BOW2")
$wsBOW2.Activate()
